# TC03 - Verify error Message for mandatory fields
#
# Duplicate the TC02 sheet (baseUrl/email test data) into a new TC03 sheet,
# change its email value to the new "00000@test.io" test input, make TC03
# the active/selected sheet, and leave TC02 selection reset back to A1.

$wb = $excel.ActiveWorkbook

$tc02 = $wb.Worksheets.Item("TC02")

# Copy TC02 and place the copy right after it -> becomes the new last sheet.
$tc02.Copy($null, $tc02)

$tc03 = $wb.Worksheets.Item($wb.Worksheets.Count)
$tc03.Name = "TC03"

# New mandatory-field test value for TC03.
$tc03.Range("B2").Value = "00000@test.io"

# Reset TC02's selection to A1 now that it is no longer the active sheet.
$tc02.Activate()
$tc02.Range("A1").Select()

# TC03 is the newly active / selected sheet, with B2 highlighted.
$tc03.Activate()
$tc03.Range("B2").Select()
